$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.108.65"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.636.72"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'214.21"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "'0.5232"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.2595"
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("D9").Value = "'0.06302"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'20.63"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").Value = "'0.07601"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.640.35"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "'4.429"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "1.862.51"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "'0.5521"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "0.0₅8007"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'65.00"
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").Value = "26.091.42"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "'4.693"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "'186.81"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'10.18"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").Value = "'6.139"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'146.01"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").Value = "'0.1214"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "'7.429"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "'15.77"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "'1.393"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("D30").Value = "'0.05911"
$ws.Range("E30").Value = "  -5.29%  "
$ws.Range("D31").Value = "'1.259"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'3.425"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "'3.402"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "'1.636"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'0.9844"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'2.390"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'2.751"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'0.5778"
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("D39").Value = "'0.01612"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'0.8551"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "1.039.66"
$ws.Range("E42").Value = "  -5.92%  "
$ws.Range("D43").Value = "'5.702"
$ws.Range("E43").Value = "  -6.90%  "
$ws.Range("D44").Value = "'100.29"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "1.788.73"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "'55.37"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.089"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "'0.9969"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "'0.05167"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'0.4220"
$ws.Range("E51").Value = "  -0.82%  "
